$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regional for Mapping")

# These "Last Cr" / "Last Date" columns (X, Y) are stored as plain text in
# this sheet even though the values look numeric/date-like, so a leading
# apostrophe is used to force literal-text entry (matching how the source
# data was typed) instead of letting it auto-convert to a number or date.

# Row 5 (R-61 S2): Last Cr / Last Date
$ws.Range("X5").Value = "'0.843"
$ws.Range("Y5").Value = "'2013-02-12"

# Row 9 (Test Well 1): Last Cr / Last Date
$ws.Range("X9").Value = "'1.5"
$ws.Range("Y9").Value = "'2005-05-04"

# Row 11 (Test Well 3): Last Cr / Last Date
$ws.Range("X11").Value = "'2.4"
$ws.Range("Y11").Value = "'2006-01-19"

# Row 12 (Test Well 4): Last Cr / Last Date
$ws.Range("X12").Value = "'2.74"
$ws.Range("Y12").Value = "'2002-05-17"

# Row 15 (R-19 S7): Last Cr only
$ws.Range("X15").Value = "'1.441"

# Row 16 (R-20 S3): Last Cr only
$ws.Range("X16").Value = "'3.7"

# Row 18 (R-22 S2): Last Cr / Last Date
$ws.Range("X18").Value = "'3.9"
$ws.Range("Y18").Value = "'2008-12-18"

# Row 23 (R-54 S1): Last Cr / Last Date
$ws.Range("X23").Value = "'1.1334"
$ws.Range("Y23").Value = "'2011-07-12"
